$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9782673716545105
$ws.Range("B1").Value = 1.760467410087585
$ws.Range("C1").Value = 4.836337566375732
$ws.Range("D1").Value = 1.270168900489807
$ws.Range("E1").Value = 1.264987230300903
